$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The price rows (2-8) are being rotated: each row's data effectively
# shifts so that the old rows 6,7,8 move up to become rows 2,3,4 and the
# old rows 2,3,4,5 move down to become rows 5,6,7,8 (cyclic rotation).
# We capture the "before" values for columns D,I,J,K,L,M,N,P,Q for rows 2-8
# and then write them back out according to the new mapping.

$cols = @("D","I","J","K","L","M","N","P","Q")

$before = @{}
for ($r = 2; $r -le 8; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowVals
}

# Mapping: new row -> old row (source of the data)
$mapping = @{
    2 = 6
    3 = 7
    4 = 8
    5 = 2
    6 = 3
    7 = 4
    8 = 5
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $src = $before[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $src[$c]
    }
}
